$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
# Row 28
$ws.Range("H28").Value = 474.7
$ws.Range("I28").Value = 476.14285
$ws.Range("J28").Value = 471.33334
$ws.Range("K28").Value = 476.14285
$ws.Range("L28").Value = 471.33334
$ws.Range("M28").Value = 8.85714999999999
$ws.Range("N28").Value = -1441.33334
# Row 55
$ws.Range("H55").Value = 230.78572
$ws.Range("I55").Value = 253.875
$ws.Range("K55").Value = 253.875
$ws.Range("M55").Value = -39.875
# Row 88
$ws.Range("H88").Value = 28976190
$ws.Range("I88").Value = 5500
$ws.Range("J88").Value = 43461536
$ws.Range("K88").Value = 5500
$ws.Range("L88").Value = 43461536
$ws.Range("M88").Value = -5094
$ws.Range("N88").Value = -43462348
# Row 91
$ws.Range("H91").Value = 28976190
$ws.Range("I91").Value = 5500
$ws.Range("J91").Value = 43461536
$ws.Range("K91").Value = 5500
$ws.Range("L91").Value = 43461536
$ws.Range("M91").Value = -4096
$ws.Range("N91").Value = -43464344

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 17048.676
$ws.Range("I32").Value = 8496.046
$ws.Range("K32").Value = 8496.046
$ws.Range("M32").Value = -8209.046
# Row 45
$ws.Range("H45").Value = 1535.0667
$ws.Range("I45").Value = 1627.1538
$ws.Range("J45").Value = 1464.6471
$ws.Range("K45").Value = 1627.1538
$ws.Range("L45").Value = 1464.6471
$ws.Range("M45").Value = -1250.1538
$ws.Range("N45").Value = -2218.6471
# Row 61
$ws.Range("H61").Value = 5183.923
$ws.Range("I61").Value = 3856.75
$ws.Range("K61").Value = 3856.75
$ws.Range("M61").Value = -3644.75
# Row 76
$ws.Range("H76").Value = 290288
$ws.Range("J76").Value = 290288
$ws.Range("L76").Value = 290288
$ws.Range("N76").Value = -290964
# Row 79
$ws.Range("H79").Value = 290288
$ws.Range("J79").Value = 290288
$ws.Range("L79").Value = 290288
$ws.Range("N79").Value = -292628
# Row 110
$ws.Range("H110").Value = 2537.5386
$ws.Range("I110").Value = 2056.5
$ws.Range("J110").Value = 3307.2
$ws.Range("K110").Value = 2056.5
$ws.Range("L110").Value = 3307.2
$ws.Range("M110").Value = -11.5
$ws.Range("N110").Value = -7397.2
# Row 136
$ws.Range("H136").Value = 5183.923
$ws.Range("I136").Value = 3856.75
$ws.Range("K136").Value = 11570.25
$ws.Range("M136").Value = -9020.25

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
# Row 75
$ws.Range("H75").Value = 11544.728
$ws.Range("I75").Value = 9688
$ws.Range("K75").Value = 9688
$ws.Range("M75").Value = -8752
# Row 78
$ws.Range("H78").Value = 11544.728
$ws.Range("I78").Value = 9688
$ws.Range("K78").Value = 29064
$ws.Range("M78").Value = -24384
# Row 140
$ws.Range("H140").Value = 90033.336

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
# Row 3
$ws.Range("H3").Value = 1000.6667
$ws.Range("I3").Value = 1000.6667
$ws.Range("K3").Value = 1000.6667
$ws.Range("M3").Value = -887.6667
# Row 31
$ws.Range("H31").Value = 433001.2
$ws.Range("I31").Value = 10067.071
$ws.Range("K31").Value = 10067.071
$ws.Range("M31").Value = -9772.071
# Row 32
$ws.Range("H32").Value = 389.5
$ws.Range("I32").Value = 389.5
$ws.Range("K32").Value = 389.5
$ws.Range("M32").Value = -73.5
# Row 34
$ws.Range("H34").Value = 433001.2
$ws.Range("I34").Value = 10067.071
$ws.Range("K34").Value = 10067.071
$ws.Range("M34").Value = -9865.071
# Row 102
$ws.Range("H102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").ClearContents()

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
# Row 40
$ws.Range("H40").Value = 140.93333
$ws.Range("J40").Value = 194.5
$ws.Range("L40").Value = 778
$ws.Range("N40").Value = -916
# Row 68
$ws.Range("H68").Value = 1741.9166
$ws.Range("J68").Value = 1907.8889
$ws.Range("L68").Value = 5723.6667
$ws.Range("N68").Value = -7345.6667
# Row 71
$ws.Range("H71").Value = 1741.9166
$ws.Range("J71").Value = 1907.8889
$ws.Range("L71").Value = 17171.0001
$ws.Range("N71").Value = -25283.0001
# Row 114
$ws.Range("H114").Value = 22677.9
$ws.Range("J114").Value = 28193.875
$ws.Range("L114").Value = 84581.625
$ws.Range("N114").Value = -91089.625
# Row 117
$ws.Range("H117").Value = 189930.83
$ws.Range("J117").Value = 227617
$ws.Range("L117").Value = 682851
$ws.Range("N117").Value = -689735
# Row 132
$ws.Range("H132").Value = 3282.5
$ws.Range("I132").Value = 1932.6666
$ws.Range("J132").Value = 4632.3335
$ws.Range("K132").Value = 17393.9994
$ws.Range("L132").Value = 41691.0015
$ws.Range("M132").Value = -14863.9994
$ws.Range("N132").Value = -46751.0015

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
# Row 11
$ws.Range("H11").Value = 3435.75
$ws.Range("I11").Value = 2999.5
$ws.Range("J11").Value = 3581.1667
$ws.Range("K11").Value = 2999.5
$ws.Range("L11").Value = 3581.1667
$ws.Range("M11").Value = -2860.5
$ws.Range("N11").Value = -3859.1667
# Row 40
$ws.Range("H40").Value = 100000
$ws.Range("J40").Value = 100000
$ws.Range("L40").Value = 100000
$ws.Range("N40").Value = -100302

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
# Row 5
$ws.Range("H5").Value = 49999
$ws.Range("I5").Value = 49999
$ws.Range("K5").Value = 49999
$ws.Range("M5").Value = -49886
# Row 17
$ws.Range("H17").Value = 10007
$ws.Range("I17").Value = 9562.5
$ws.Range("J17").Value = 10599.667
$ws.Range("K17").Value = 9562.5
$ws.Range("L17").Value = 10599.667
$ws.Range("N17").Value = -10939.667
$ws.Range("M17").Value = -9392.5
# Row 132
$ws.Range("H132").Value = 7516.7354
$ws.Range("I132").Value = 7502.6553
$ws.Range("J132").Value = 7598.4
$ws.Range("K132").Value = 22507.9659
$ws.Range("L132").Value = 22795.2
$ws.Range("M132").Value = -19977.9659
$ws.Range("N132").Value = -27855.2

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
# Row 96
$ws.Range("H96").Value = 1325360.2
$ws.Range("I96").Value = 3089937
$ws.Range("J96").Value = 1927.6875
$ws.Range("K96").Value = 3089937
$ws.Range("L96").Value = 1927.6875
$ws.Range("M96").Value = -3088564
$ws.Range("N96").Value = -4673.6875
# Row 113
$ws.Range("H113").Value = 501.58334
$ws.Range("I113").Value = 483.85715
$ws.Range("K113").Value = 1451.57145
$ws.Range("M113").Value = 718.4285500000001
# Row 122
$ws.Range("H122").Value = 26392604
$ws.Range("I122").Value = 25003304
$ws.Range("K122").Value = 75009912
$ws.Range("M122").Value = -75007462
